# ============================================================
# feat: add 2022-Q1 data
# - insert a new "2022-Q1" worksheet (between "2021-Q4" and "总计")
#   holding the per-fund holding breakdown for that quarter
# - add the matching summary row to the "总计" (totals) sheet
# ============================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before "总计"
# ---------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q4")   # same column layout/style to copy from
$totalSheet = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# Header row (row 1) -- columns B:H -- same text/style as the other
# quarterly sheets (bold, centered, thin border), so copy it verbatim
$refSheet.Range("B1:H1").Copy($ws.Range("B1:H1"))

# Columns B:G hold text values (fund code / name / formatted decimal
# strings) in the source data -- force text storage up front so Excel
# does not auto-convert numeric-looking strings (e.g. "001445" -> 1445)
$ws.Range("B2:G20").NumberFormat = "@"

# Data rows 2-20 (index column A holds 0..18, styled like the other sheets A column)

# row 2 (index 0)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(2,1))
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "001445"
$ws.Cells.Item(2,3).Value = "华安国企改革主题灵活配置混合"
$ws.Cells.Item(2,4).Value = "47.50"
$ws.Cells.Item(2,5).Value = "87.35"
$ws.Cells.Item(2,6).Value = "2.73"
$ws.Cells.Item(2,7).Value = "1.2968"
$ws.Cells.Item(2,8).Value = 10

# row 3 (index 1)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(3,1))
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "006682"
$ws.Cells.Item(3,3).Value = "景顺长城中证500指数增强"
$ws.Cells.Item(3,4).Value = "16.63"
$ws.Cells.Item(3,5).Value = "87.75"
$ws.Cells.Item(3,6).Value = "2.05"
$ws.Cells.Item(3,7).Value = "0.3409"
$ws.Cells.Item(3,8).Value = 3

# row 4 (index 2)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(4,1))
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "501029"
$ws.Cells.Item(4,3).Value = "华宝标普中国A股红利机会指数（LOF）A"
$ws.Cells.Item(4,4).Value = "13.19"
$ws.Cells.Item(4,5).Value = "94.39"
$ws.Cells.Item(4,6).Value = "2.23"
$ws.Cells.Item(4,7).Value = "0.2941"
$ws.Cells.Item(4,8).Value = 2

# row 5 (index 3)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(5,1))
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "080005"
$ws.Cells.Item(5,3).Value = "长盛量化红利混合"
$ws.Cells.Item(5,4).Value = "2.66"
$ws.Cells.Item(5,5).Value = "69.88"
$ws.Cells.Item(5,6).Value = "3.16"
$ws.Cells.Item(5,7).Value = "0.0841"
$ws.Cells.Item(5,8).Value = 4

# row 6 (index 4)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(6,1))
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "008851"
$ws.Cells.Item(6,3).Value = "景顺长城量化对冲策略三个月定期开放灵活配置混合"
$ws.Cells.Item(6,4).Value = "5.05"
$ws.Cells.Item(6,5).Value = "74.55"
$ws.Cells.Item(6,6).Value = "1.65"
$ws.Cells.Item(6,7).Value = "0.0833"
$ws.Cells.Item(6,8).Value = 1

# row 7 (index 5)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(7,1))
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "010857"
$ws.Cells.Item(7,3).Value = "宝盈祥乐一年持有期混合型证券投资基金A"
$ws.Cells.Item(7,4).Value = "2.00"
$ws.Cells.Item(7,5).Value = "36.12"
$ws.Cells.Item(7,6).Value = "3.61"
$ws.Cells.Item(7,7).Value = "0.0722"
$ws.Cells.Item(7,8).Value = 3

# row 8 (index 6)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(8,1))
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "004945"
$ws.Cells.Item(8,3).Value = "长信中证500指数增强"
$ws.Cells.Item(8,4).Value = "2.81"
$ws.Cells.Item(8,5).Value = "92.77"
$ws.Cells.Item(8,6).Value = "1.64"
$ws.Cells.Item(8,7).Value = "0.0461"
$ws.Cells.Item(8,8).Value = 10

# row 9 (index 7)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(9,1))
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "008324"
$ws.Cells.Item(9,3).Value = "宝盈祥利稳健配置混合A"
$ws.Cells.Item(9,4).Value = "1.09"
$ws.Cells.Item(9,5).Value = "36.53"
$ws.Cells.Item(9,6).Value = "3.96"
$ws.Cells.Item(9,7).Value = "0.0432"
$ws.Cells.Item(9,8).Value = 3

# row 10 (index 8)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(10,1))
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "004258"
$ws.Cells.Item(10,3).Value = "国寿安保稳嘉混合A"
$ws.Cells.Item(10,4).Value = "2.53"
$ws.Cells.Item(10,5).Value = "22.03"
$ws.Cells.Item(10,6).Value = "0.71"
$ws.Cells.Item(10,7).Value = "0.0180"
$ws.Cells.Item(10,8).Value = 10

# row 11 (index 9)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(11,1))
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "004301"
$ws.Cells.Item(11,3).Value = "国寿安保稳信混合A"
$ws.Cells.Item(11,4).Value = "1.50"
$ws.Cells.Item(11,5).Value = "20.03"
$ws.Cells.Item(11,6).Value = "1.20"
$ws.Cells.Item(11,7).Value = "0.0180"
$ws.Cells.Item(11,8).Value = 6

# row 12 (index 10)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(12,1))
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "008325"
$ws.Cells.Item(12,3).Value = "宝盈祥利稳健配置混合C"
$ws.Cells.Item(12,4).Value = "0.43"
$ws.Cells.Item(12,5).Value = "36.53"
$ws.Cells.Item(12,6).Value = "3.96"
$ws.Cells.Item(12,7).Value = "0.0170"
$ws.Cells.Item(12,8).Value = 3

# row 13 (index 11)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(13,1))
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "512590"
$ws.Cells.Item(13,3).Value = "浦银安盛中证高股息精选ETF"
$ws.Cells.Item(13,4).Value = "0.59"
$ws.Cells.Item(13,5).Value = "96.43"
$ws.Cells.Item(13,6).Value = "2.41"
$ws.Cells.Item(13,7).Value = "0.0142"
$ws.Cells.Item(13,8).Value = 5

# row 14 (index 12)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(14,1))
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "008112"
$ws.Cells.Item(14,3).Value = "中泰中证500指数增强A"
$ws.Cells.Item(14,4).Value = "0.61"
$ws.Cells.Item(14,5).Value = "92.46"
$ws.Cells.Item(14,6).Value = "1.64"
$ws.Cells.Item(14,7).Value = "0.0100"
$ws.Cells.Item(14,8).Value = 3

# row 15 (index 13)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(15,1))
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "008113"
$ws.Cells.Item(15,3).Value = "中泰中证500指数增强C"
$ws.Cells.Item(15,4).Value = "0.46"
$ws.Cells.Item(15,5).Value = "92.46"
$ws.Cells.Item(15,6).Value = "1.64"
$ws.Cells.Item(15,7).Value = "0.0075"
$ws.Cells.Item(15,8).Value = 3

# row 16 (index 14)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(16,1))
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "010858"
$ws.Cells.Item(16,3).Value = "宝盈祥乐一年持有期混合型证券投资基金C"
$ws.Cells.Item(16,4).Value = "0.14"
$ws.Cells.Item(16,5).Value = "36.12"
$ws.Cells.Item(16,6).Value = "3.61"
$ws.Cells.Item(16,7).Value = "0.0051"
$ws.Cells.Item(16,8).Value = 3

# row 17 (index 15)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(17,1))
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "005770"
$ws.Cells.Item(17,3).Value = "信达澳银中证沪港深高股息精选指数"
$ws.Cells.Item(17,4).Value = "0.01"
$ws.Cells.Item(17,5).Value = "92.47"
$ws.Cells.Item(17,6).Value = "2.17"
$ws.Cells.Item(17,7).Value = "0.0002"
$ws.Cells.Item(17,8).Value = 9

# row 18 (index 16)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(18,1))
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "004302"
$ws.Cells.Item(18,3).Value = "国寿安保稳信混合C"
$ws.Cells.Item(18,4).Value = "0.01"
$ws.Cells.Item(18,5).Value = "20.03"
$ws.Cells.Item(18,6).Value = "1.20"
$ws.Cells.Item(18,7).Value = "0.0001"
$ws.Cells.Item(18,8).Value = 6

# row 19 (index 17)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(19,1))
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "004259"
$ws.Cells.Item(19,3).Value = "国寿安保稳嘉混合C"
$ws.Cells.Item(19,4).Value = "0.00"
$ws.Cells.Item(19,5).Value = "22.03"
$ws.Cells.Item(19,6).Value = "0.71"
$g = $ws.Cells.Item(19,7)
$g.NumberFormat = "General"
$g.Value = 0
$ws.Cells.Item(19,8).Value = 10

# row 20 (index 18)
$refSheet.Cells.Item(2,1).Copy($ws.Cells.Item(20,1))
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "015406"
$ws.Cells.Item(20,3).Value = "国寿安保稳信混合E"
$ws.Cells.Item(20,4).Value = "0.00"
$ws.Cells.Item(20,5).Value = "20.03"
$ws.Cells.Item(20,6).Value = "1.20"
$g = $ws.Cells.Item(20,7)
$g.NumberFormat = "General"
$g.Value = 0
$ws.Cells.Item(20,8).Value = 6

# ---------------------------------------------------------------
# 2. Add the 2022-Q1 summary row to the "总计" sheet
# ---------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()

# New row 2 reuses the existing index-column style from row 3 (the old row 2,
# pushed down by the insert) so it matches the other rows exactly
$zj.Cells.Item(3,1).Copy($zj.Cells.Item(2,1))
$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = "2022-Q1"
$zj.Cells.Item(2,3).Value = 19
$zj.Cells.Item(2,4).Value = 2.35

# Renumber the index column (A) for the rows pushed down by the insert
# (old index 0..4 in rows 2..6 is now rows 3..7, index should be 1..5)
for ($row = 3; $row -le 7; $row++) {
    $zj.Cells.Item($row, 1).Value = $row - 2
}
